$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.852.02"
$ws.Range("E2").Value = "  +2.26%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.874.80"
$ws.Range("E3").Value = "  +0.77%  "

# Row 4
$ws.Range("E4").Value = "  -0.71%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.50"
$ws.Range("E5").Value = "  +0.52%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.012"
$ws.Range("E6").Value = "  -0.69%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4829"
$ws.Range("E7").Value = "  +0.65%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3811"
$ws.Range("E8").Value = "  +2.31%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07374"
$ws.Range("E9").Value = "  +0.89%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9400"
$ws.Range("E10").Value = "  +0.34%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.00"
$ws.Range("E11").Value = "  +3.55%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07785"
$ws.Range("E12").Value = "  -1.13%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.889.63"
$ws.Range("E13").Value = "  +1.92%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.522"
$ws.Range("E14").Value = "  +1.82%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.600"
$ws.Range("E15").Value = "  +1.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.29"
$ws.Range("E16").Value = "  +1.10%  "

# Row 17
$ws.Range("E17").Value = "  -0.78%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008868"
$ws.Range("E18").Value = "  +1.41%  "

# Row 19
$ws.Range("E19").Value = "  -0.73%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.875.96"
$ws.Range("E20").Value = "  +2.26%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.85"
$ws.Range("E21").Value = "  +0.96%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.123"
$ws.Range("E22").Value = "  +0.36%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.120.93"
$ws.Range("E23").Value = "  +1.88%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.92"
$ws.Range("E24").Value = "  +2.41%  "

# Row 25
$ws.Range("E25").Value = "  +0.01%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.56"
$ws.Range("E26").Value = "  +2.41%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.55"
$ws.Range("E27").Value = "  +0.32%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.041"
$ws.Range("E28").Value = "  +2.18%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.98"
$ws.Range("E29").Value = "  +0.29%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.968"
$ws.Range("E30").Value = "  -0.13%  "

# Row 31
$ws.Range("E31").Value = "  -0.06%  "

# Row 32
$ws.Range("E32").Value = "  -0.16%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.224"
$ws.Range("E33").Value = "  +3.39%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7728"
$ws.Range("E34").Value = "  +4.14%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.654"
$ws.Range("E35").Value = "  +1.50%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.729"
$ws.Range("E36").Value = "  +1.55%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02046"
$ws.Range("E37").Value = "  +1.01%  "

# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.125"
$ws.Range("E38").Value = "  +0.16%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5607"
$ws.Range("E39").Value = "  +5.20%  "

# Row 40
$ws.Range("E40").Value = "  +2.19%  "

# Row 41
$ws.Range("E41").Value = "  +0.07%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.047"
$ws.Range("E42").Value = "  -0.86%  "

# Row 43
$ws.Range("E43").Value = "  +2.48%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1529"
$ws.Range("E44").Value = "  +0.08%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.71"
$ws.Range("E45").Value = "  +1.11%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4877"
$ws.Range("E46").Value = "  +1.91%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.62"
$ws.Range("E47").Value = "  +2.83%  "

# Row 48
$ws.Range("E48").Value = "  -0.71%  "

# Row 49
$ws.Range("E49").Value = "  +1.90%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.13"
$ws.Range("E50").Value = "  +2.70%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06117"
$ws.Range("E51").Value = "  +0.69%  "
